$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column A, shifting all other columns (B:F) one position to the left.
$ws.Range("A:A").Delete()
